$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overordnet projektplan")

# --- Copy formats first, while source cells still have their original style ---
# D21 needs style s="8" (same as D20's current style) but stays empty.
$ws.Range("D20").Copy()
$ws.Range("D21").PasteSpecial(-4122)

# D16 needs style s="7" (italic, no alignment) - same style as E15/F15/G15 currently have.
$ws.Range("E15").Copy()
$ws.Range("D16").PasteSpecial(-4122)

# E9, E10, E11 need style s="20" (same style E8 currently has).
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 6: add new cell K6, keep AC6 value the same ("Aflevering af projekt") ---
$ws.Range("K6").Value = "Review af artefakter med gruppe grp 2"
$ws.Range("AC6").Value = "Aflevering af projekt"

# --- Row 7: fill E7 ---
$ws.Range("E7").Value = "1 time"

# --- Row 8: fill E8 ---
$ws.Range("E8").Value = "0,5 timer"

# --- Row 9: add E9 ---
$ws.Range("E9").Value = "0,5 timer"

# --- Row 10: change D10, add E10 ---
$ws.Range("D10").Value = "Lav UC1"
$ws.Range("E10").Value = "1 time"

# --- Row 11: change D11, add E11 ---
$ws.Range("D11").Value = "Lav UC2"
$ws.Range("E11").Value = "0,5 timer"

# --- Row 12: change D12 ---
$ws.Range("D12").Value = "Mockups"

# --- Row 13: add D13, E13 ---
$ws.Range("D13").Value = "MUST-undersøgelse"
$ws.Range("E13").Value = "1 time"

# --- Row 15: remove D15 entirely (cell element disappears) ---
$ws.Range("D15").Clear()

# --- Row 16: clear D16 content (style already set above to s=7) ---
$ws.Range("D16").ClearContents()

# --- Rows 17-20: clear D content but keep style s=8 ---
$ws.Range("D17").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("D19").ClearContents()
$ws.Range("D20").ClearContents()

# --- Selection change recorded in the sheet view ---
$ws.Range("E12").Select()
